$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (search_radius=4 / #PCs=5) entirely, which
# shrinks the used range from A1:H6 down to A1:H5.
$ws.Rows.Item(6).Delete()

# Row 2 (search_radius bucket 0): search radius shrank 120000 -> 90000,
# so the variogram-derived statistics were recomputed.
$ws.Range("B2").Value = 90000
$ws.Range("D2").Value = 0.5379465951389433
$ws.Range("E2").Value = 3.296988220851615
$ws.Range("F2").Value = 1.006
$ws.Range("H2").Value = 3.68572083840303

# Row 3 (search_radius bucket 1)
$ws.Range("B3").Value = 90000
$ws.Range("D3").Value = 0.7281377488446686
$ws.Range("E3").Value = 3.150579992672729
$ws.Range("F3").Value = 1.453
$ws.Range("H3").Value = 3.68572083840303

# Row 4 (search_radius bucket 2)
$ws.Range("B4").Value = 90000
$ws.Range("D4").Value = 0.8070082382439561
$ws.Range("E4").Value = 3.14240151487407
$ws.Range("F4").Value = 1.7876
$ws.Range("H4").Value = 3.68572083840303

# Row 5 (search_radius bucket 3)
$ws.Range("B5").Value = 90000
$ws.Range("D5").Value = 0.9689016246478028
$ws.Range("E5").Value = 3.060298046692686
$ws.Range("F5").Value = 1.9857
$ws.Range("H5").Value = 3.68572083840303
